$d = $word.ActiveDocument

# 1. "Template" -> "Essential Elements"
$d.Content.Find.Execute("Template", $false, $false, $false, $false, $false, $true, 1, $false, "Essential Elements", 2)

# 2. Remove the 4 blank paragraphs following "Post a picture..." paragraph.
#    (Originally paragraphs 19-22, right before "Write a short paragraph...")
$pPost = $d.Paragraphs.Item(18)
$pWrite = $d.Paragraphs.Item(23)
$rBlank2 = $d.Range($pPost.Range.End, $pWrite.Range.Start)
$rBlank2.Delete()

# 3. Remove the old "Provide the date, time and brief description of the activity."
#    paragraph's 4 trailing blank paragraphs, and the paragraph itself; we keep the
#    paragraph's own (numbered) formatting by writing the combined sentence into it,
#    then delete the now-redundant empty bullet paragraph that preceded it.
$pOldProvide = $d.Paragraphs.Item(13)
$pPost2 = $d.Paragraphs.Item(18)
$rBlank1 = $d.Range($pOldProvide.Range.End, $pPost2.Range.Start)
$rBlank1.Delete()

# Set the surviving "Provide the date, time..." paragraph's text as three runs.
$pOldProvide.Range.Text = "Provide the date, time"
$pOldProvide.Range.InsertAfter(",")
$pOldProvide.Range.InsertAfter(" and brief description of the activity. ")

# 4. Delete the now-obsolete empty bullet paragraph (the one right after
#    "Include your name:") that used to separate the two items.
$pEmpty = $d.Paragraphs.Item(12)
$pEmpty.Range.Delete()

Write-Host "Done"
